$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    current "2022-Q2" sheet), pre-populated with the per-fund holdings
#    table.  We clone "2021-Q2" as a template since it already has the
#    exact row count (1 header row + 6 data rows) and cell styling we
#    need, then simply overwrite every cell's content.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q2")
$anchor   = $wb.Worksheets.Item("2022-Q2")
$template.Copy($anchor)
$q3 = $wb.Worksheets.Item("2021-Q2 (2)")
$q3.Name = "2022-Q3"

# Header row
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# Row 2
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "513090"
$q3.Cells.Item(2,3).Value = "易方达中证香港证券投资主题ETF"
$q3.Cells.Item(2,4).Value = "10.53"
$q3.Cells.Item(2,5).Value = "96.33"
$q3.Cells.Item(2,6).Value = "8.23"
$q3.Cells.Item(2,7).Value = "0.8666"
$q3.Cells.Item(2,8).Value = 5

# Row 3
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "004497"
$q3.Cells.Item(3,3).Value = "前海开源多元策略灵活配置混合C"
$q3.Cells.Item(3,4).Value = "1.68"
$q3.Cells.Item(3,5).Value = "93.04"
$q3.Cells.Item(3,6).Value = "4.39"
$q3.Cells.Item(3,7).Value = "0.0738"
$q3.Cells.Item(3,8).Value = 7

# Row 4
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).Value = "004496"
$q3.Cells.Item(4,3).Value = "前海开源多元策略灵活配置混合A"
$q3.Cells.Item(4,4).Value = "0.91"
$q3.Cells.Item(4,5).Value = "93.04"
$q3.Cells.Item(4,6).Value = "4.39"
$q3.Cells.Item(4,7).Value = "0.0399"
$q3.Cells.Item(4,8).Value = 7

# Row 5
$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).Value = "006106"
$q3.Cells.Item(5,3).Value = "景顺长城量化港股通股票"
$q3.Cells.Item(5,4).Value = "0.52"
$q3.Cells.Item(5,5).Value = "50.20"
$q3.Cells.Item(5,6).Value = "0.95"
$q3.Cells.Item(5,7).Value = "0.0049"
$q3.Cells.Item(5,8).Value = 10

# Row 6
$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).Value = "004532"
$q3.Cells.Item(6,3).Value = "民生加银中证港股通高股息精选指数A"
$q3.Cells.Item(6,4).Value = "0.13"
$q3.Cells.Item(6,5).Value = "92.87"
$q3.Cells.Item(6,6).Value = "3.68"
$q3.Cells.Item(6,7).Value = "0.0048"
$q3.Cells.Item(6,8).Value = 8

# Row 7
$q3.Cells.Item(7,1).Value = 5
$q3.Cells.Item(7,2).Value = "004533"
$q3.Cells.Item(7,3).Value = "民生加银中证港股通高股息精选指数C"
$q3.Cells.Item(7,4).Value = "0.08"
$q3.Cells.Item(7,5).Value = "92.87"
$q3.Cells.Item(7,6).Value = "3.68"
$q3.Cells.Item(7,7).Value = "0.0029"
$q3.Cells.Item(7,8).Value = 8

# ---------------------------------------------------------------------
# 2) Add the corresponding summary row on "总计": insert a blank row
#    right under the header and fill it with the 2022-Q3 totals; then
#    renumber the sequential index in column A for every row below it.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Copy the column-A style (bold + border + centered) from the row below
# onto the freshly inserted A2 cell.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 0.99

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(8,1).Value = 6
$total.Cells.Item(9,1).Value = 7

# Restore "总计" as the active sheet (sheet insertion/copy operations
# above moved the selection onto the newly created sheet).
$total.Activate()
